# Wealth Allocation: border styling
#
# The "Financial Statement" sheet had a stray unused row (row 9) sitting
# between the "Income & Expenses" block and the "Wealth Allocation" block.
# Removing that entire row shifts everything below it up by one, which is
# exactly what happened in the authored edit (old row 10 -> new row 9,
# old row 11 "Wealth Allocation" header -> new row 10, etc).

$wb = $excel.ActiveWorkbook

$wsFin = $wb.Worksheets.Item("Financial Statement")
$wsFin.Rows.Item(9).Delete()

# AutoFit column A now that the layout settled (matches the <cols> width
# that Excel persists after an autofit pass).
$wsFin.Columns.Item(1).AutoFit()

# Restore/refresh the selections Excel remembers per sheet, and make sure
# "Financial Statement" is the active tab (it was previously "Records").
$wsFin.Activate()
$wsFin.Range("O8").Select()

$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Activate()
$wsAssets.Range("A22").Select()

$wsExpenses = $wb.Worksheets.Item("Expenses")
$wsExpenses.Activate()
$wsExpenses.Range("B21").Select()

$wsFin.Activate()
